# Updates the cryptos list (Price + Volume(1h) columns) for rows 2-51.
# Price values (column D) are forced to remain plain text (matching the
# original inline-string cells) via a leading quote-prefix, then the style
# is reset to "Normal" so no stray number-format/style is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''66.825.89'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +8.10%  '
$ws.Range("D3").Value = '''3.866.10'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +12.60%  '
$ws.Range("D4").Value = '''0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -1.23%  '
$ws.Range("D5").Value = '''424.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +9.25%  '
$ws.Range("D6").Value = '''131.65'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.66%  '
$ws.Range("D7").Value = '''3.856.86'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +12.10%  '
$ws.Range("D8").Value = '''0.614'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.74%  '
$ws.Range("D9").Value = '''0.998'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("D10").Value = '''0.730'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.80%  '
$ws.Range("E11").Value = '  +11.34%  '
$ws.Range("E12").Value = '  +23.60%  '
$ws.Range("D13").Value = '''41.01'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.16%  '
$ws.Range("D14").Value = '''10.27'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +11.49%  '
$ws.Range("D15").Value = '''4.465.84'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +9.82%  '
$ws.Range("D16").Value = '''15.93'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +28.20%  '
$ws.Range("D17").Value = '''3.864.63'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +11.20%  '
$ws.Range("E18").Value = '  +0.44%  '
$ws.Range("D19").Value = '''19.97'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.32%  '
$ws.Range("D20").Value = '''66.919.55'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.61%  '
$ws.Range("E21").Value = '  +6.80%  '
$ws.Range("D22").Value = '''413.43'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.87%  '
$ws.Range("D23").Value = '''14.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.93%  '
$ws.Range("D24").Value = '''84.48'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.94%  '
$ws.Range("D25").Value = '''3.03'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.03%  '
$ws.Range("D26").Value = '''37.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +13.89%  '
$ws.Range("D27").Value = '''9.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +12.74%  '
$ws.Range("D28").Value = '''3.26'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.42%  '
$ws.Range("D29").Value = '''5.25'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.29%  '
$ws.Range("D30").Value = '''9.09'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +39.46%  '
$ws.Range("D31").Value = '''725.62'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +13.12%  '
$ws.Range("D32").Value = '''13.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +11.95%  '
$ws.Range("E33").Value = '  +11.93%  '
$ws.Range("E34").Value = '  +6.09%  '
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("D36").Value = '''39.12'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.06%  '
$ws.Range("D37").Value = '''0.153'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.02%  '
$ws.Range("D38").Value = '''55.56'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.93%  '
$ws.Range("D39").Value = '''5.45'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +34.80%  '
$ws.Range("D40").Value = '''0.0{0}0758' -f [char]0x2083
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +34.57%  '
$ws.Range("D41").Value = '''0.0462'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.59%  '
$ws.Range("E42").Value = '  +7.63%  '
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("E44").Value = '  +1.88%  '
$ws.Range("E45").Value = '  +8.87%  '
$ws.Range("D46").Value = '''3.14'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.44%  '
$ws.Range("D47").Value = '''0.314'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +16.33%  '
$ws.Range("E48").Value = '  +6.45%  '
$ws.Range("D49").Value = '''141.54'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("E50").Value = '  +6.16%  '
$ws.Range("D51").Value = '''2.57'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.70%  '
